# Update column G ("K") values on Sheet1 for rows 2-17.
# These values were regenerated to count "K" occurrences instead of "Strike#".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 7
    3  = 2
    4  = 4
    5  = 0
    6  = 5
    7  = 4
    8  = 10
    9  = 12
    10 = 5
    11 = 5
    12 = 3
    13 = 4
    14 = 7
    15 = 2
    16 = 6
    17 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
